$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), matching the style of the
# existing header cell H1 ("IP") by copying its formatting.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for new columns I (I0) and J (IF), rows 2-24
$i0 = @(8,8,4,1,1,6,1,1,1,1,1,1,1,2,1,1,1,1,1,1,5,4,3)
$if = @(9,8,5,5,5,8,5,3,2,6,4,4,6,4,4,4,4,2,4,5,7,5,3)

for ($idx = 0; $idx -lt $i0.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $i0[$idx]
    $ws.Cells.Item($row, 10).Value = $if[$idx]
}
